$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.357.51'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.87%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.425.70'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.26%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '512.16'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.13'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.26%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.548'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.435.65'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.88%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0952'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -5.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.17'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.332'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.63%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.855.69'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '57.256.15'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.72'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.02%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.431.23'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.43'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '314.84'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.10'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.66'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.37%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.53'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.16%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.82%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.159'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.19'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.45%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '168.60'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.14%  '

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0720'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.11%  '

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.22'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.66'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.55%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.16'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.86%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.67'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.27'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.89'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.29%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.11'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.45'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.90%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.776'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.04%  '

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.36'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.42%  '

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.97'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.66%  '

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '267.84'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.585'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.81%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0905'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.08%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '120.09'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0483'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.03'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0209'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.43'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.86%  '
